# "added 4wk low sales check" -- refreshed forecast figures on the
# "Forecast Comparison" sheet (MyForecast / Inventory Coverage / Seasonality
# Index) together with the dependent roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H), Seasonality Index (L) ---
$wsForecast.Range("D2").Value = 23
$wsForecast.Range("H2").Value = 22.09
$wsForecast.Range("L2").Value = 0.98

$wsForecast.Range("D3").Value = 23
$wsForecast.Range("H3").Value = 21.09
$wsForecast.Range("L3").Value = 0.93

$wsForecast.Range("D4").Value = 23
$wsForecast.Range("H4").Value = 20.09
$wsForecast.Range("L4").Value = 0.88

$wsForecast.Range("D5").Value = 23
$wsForecast.Range("H5").Value = 19.09
$wsForecast.Range("L5").Value = 1.05

$wsForecast.Range("D6").Value = 23
$wsForecast.Range("H6").Value = 18.09
$wsForecast.Range("L6").Value = 0.97

$wsForecast.Range("D7").Value = 23
$wsForecast.Range("H7").Value = 17.09
$wsForecast.Range("L7").Value = 0.91

$wsForecast.Range("D8").Value = 23
$wsForecast.Range("H8").Value = 16.09
$wsForecast.Range("L8").Value = 1.08

$wsForecast.Range("D9").Value = 23
$wsForecast.Range("H9").Value = 15.09
$wsForecast.Range("L9").Value = 1.15

$wsForecast.Range("D10").Value = 23
$wsForecast.Range("H10").Value = 14.09
$wsForecast.Range("L10").Value = 1.12

$wsForecast.Range("D11").Value = 23
$wsForecast.Range("H11").Value = 13.09
$wsForecast.Range("L11").Value = 1.03

$wsForecast.Range("D12").Value = 23
$wsForecast.Range("H12").Value = 12.09
$wsForecast.Range("L12").Value = 1.05

$wsForecast.Range("D13").Value = 22
$wsForecast.Range("H13").Value = 11.59
$wsForecast.Range("L13").Value = 1.19

$wsForecast.Range("D14").Value = 23
$wsForecast.Range("H14").Value = 10.13
$wsForecast.Range("L14").Value = 0.86

$wsForecast.Range("D15").Value = 23
$wsForecast.Range("H15").Value = 9.13
$wsForecast.Range("L15").Value = 0.83

$wsForecast.Range("D16").Value = 23
$wsForecast.Range("H16").Value = 8.13

$wsForecast.Range("D17").Value = 23
$wsForecast.Range("H17").Value = 7.13
$wsForecast.Range("L17").Value = 0.88

# --- Summary: forecast roll-up totals (stored as text on this sheet) ---
$wsSummary.Range("B9").Value = "367"
$wsSummary.Range("B10").Value = "184"
$wsSummary.Range("B11").Value = "92"
$wsSummary.Range("B12").Value = "23"
$wsSummary.Range("B14").Value = "22"
